$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3
$ws.Range('C3').Value = 0.216494553753843
$ws.Range('D3').Value = 0.9623306322514821
$ws.Range('E3').Value = 0.3764006536953189
$ws.Range('G3').Value = 'max\_depth: 2, max\_features: 4 \\'

# Row 4
$ws.Range('C4').Value = 0.2169259982527708
$ws.Range('D4').Value = 0.9642484276428003
$ws.Range('E4').Value = 0.3691163795579652
$ws.Range('G4').Value = 'max\_depth: 4, max\_features: 6, n\_estimators: 200 \\'

# Row 5
$ws.Range('C5').Value = 0.2071409276599118
$ws.Range('D5').Value = 0.9207532310802201
$ws.Range('E5').Value = 0.353826644576964

# Row 6
$ws.Range('C6').Value = 0.2066927184661046
$ws.Range('D6').Value = 0.9187609156645271
$ws.Range('E6').Value = 0.3550930811283948

# Row 7
$ws.Range('C7').Value = 0.213517639720418
$ws.Range('D7').Value = 0.9490980796801988
$ws.Range('E7').Value = 0.3693708632465984
$ws.Range('G7').Value = 'learning\_rate: 0.1, max\_depth: 2, max\_features: 2, n\_estimators: 50 \\'

# Row 10
$ws.Range('C10').Value = 0.217155695367041
$ws.Range('D10').Value = 0.9652694444091308
$ws.Range('E10').Value = 0.3675045227594317
$ws.Range('G10').Value = 'n\_hidden: 2, n\_neurons: 70, activation: gelu, learning\_rate: 0.007, input\_shape: 12 \\'

# Row 11
$ws.Range('C11').Value = 0.2428509513926233
$ws.Range('D11').Value = 1.079486322146746
$ws.Range('E11').Value = 0.4389462950032964

# Row 12
$ws.Range('C12').Value = 0.2079875961148583
$ws.Range('D12').Value = 0.9245167206250081
$ws.Range('E12').Value = 0.3649943448404027

# Row 13
$ws.Range('C13').Value = 0.2076820698894318
$ws.Range('D13').Value = 0.9231586391371087
$ws.Range('E13').Value = 0.3558539075391847

# Row 14
$ws.Range('C14').Value = 0.238694385469575
$ws.Range('D14').Value = 1.061010149682519
$ws.Range('E14').Value = 0.3926799918824795

# Row 15
$ws.Range('C15').Value = 0.232172000198546
$ws.Range('D15').Value = 1.032017775357972
$ws.Range('E15').Value = 0.3933912945452072

# Row 16
$ws.Range('C16').Value = 0.2048793892832957
$ws.Range('D16').Value = 0.9107005640819343
$ws.Range('E16').Value = 0.3520750957577883
$ws.Range('F16').Value = 1
$ws.Range('G16').Value = 'omega: 10000, r: 0.1 \\'

# Row 17
$ws.Range('C17').Value = 0.2048778744629491
$ws.Range('D17').Value = 0.9106938306191446
$ws.Range('E17').Value = 0.3520657595888778
$ws.Range('F17').Value = 1
$ws.Range('G17').Value = 'omega: 1000, r: 0.1 \\'

# Row 18
$ws.Range('C18').Value = 0.2040855971238866
$ws.Range('D18').Value = 0.9071721126849119
$ws.Range('E18').Value = 0.3521733687739505
$ws.Range('F18').Value = 2
$ws.Range('G18').Value = 'mu: 0.3, omega: 1000 \\'

# Row 19
$ws.Range('C19').Value = 0.2054998676977693
$ws.Range('D19').Value = 0.9134586260033333
$ws.Range('E19').Value = 0.3530993065361773
$ws.Range('F19').Value = 1
$ws.Range('G19').Value = 'alpha: 0.1, beta: 0.1, lambda1: 0.001, omega: 100, sigma: 0.1 \\'

# Row 20
$ws.Range('C20').Value = 0.4508895687924844
$ws.Range('D20').Value = 2.004229835291953
$ws.Range('E20').Value = 0.6179100142388599
$ws.Range('F20').Value = 138

# Row 21
$ws.Range('C21').Value = 0.205449459453054
$ws.Range('D21').Value = 0.9132345584821582
$ws.Range('E21').Value = 0.3532152811008513
$ws.Range('F21').Value = 1
$ws.Range('G21').Value = 'alpha: 0.001, beta: 0.01, e\_utility: 0.05, lambda1: 0.5, omega: 10000, pi: 0.5, sigma: 0.5 \\'

# Row 22
$ws.Range('C22').Value = 0.2316480938599464
$ws.Range('D22').Value = 1.02968898181872
$ws.Range('E22').Value = 0.3794632915924392
$ws.Range('F22').Value = 20
$ws.Range('G22').Value = 'alpha: 0.05, beta: 0.25, e\_utility: 0.05, lambda1: 0.001, sigma: 10 \\'

# Row 23
$ws.Range('C23').Value = 0.2848126221400735
$ws.Range('D23').Value = 1.266008340555745
$ws.Range('E23').Value = 0.4084389984485181
$ws.Range('F23').Value = 13
$ws.Range('G23').Value = 'fuzzy\_operator: prod, rules: 13 \\'

# Row 24
$ws.Range('C24').Value = 0.2371424161043116
$ws.Range('D24').Value = 1.054111557387181
$ws.Range('E24').Value = 0.4149328391550006

# Row 25
$ws.Range('C25').Value = 0.2103253550790418
$ws.Range('D25').Value = 0.934908192479825
$ws.Range('E25').Value = 0.3552525209252443
$ws.Range('F25').Value = 2
$ws.Range('G25').Value = 'adaptive\_filter: wRLS, fuzzy\_operator: prod, rules: 2 \\'

# Row 26
$ws.Range('C26').Value = 0.20731282936038
$ws.Range('D26').Value = 0.9215173439376954
$ws.Range('E26').Value = 0.3535348594976804
$ws.Range('F26').Value = 11
$ws.Range('G26').Value = 'error\_metric: RMSE, fuzzy\_operator: minmax, num\_generations: 10, num\_parents\_mating: 5, parallel\_processing: 10, rules: 11, sol\_per\_pop: 10 \\'

# Row 27
$ws.Range('C27').Value = 0.2090063323271945
$ws.Range('D27').Value = 0.9290450611597525
$ws.Range('E27').Value = 0.3518088460135911
$ws.Range('G27').Value = 'adaptive\_filter: RLS, error\_metric: MAE, fuzzy\_operator: prod, lambda1: 0.97, num\_generations: 5, num\_parents\_mating: 5, parallel\_processing: 10, rules: 1, sol\_per\_pop: 5 \\'

# Row 28
$ws.Range('C28').Value = 0.2404912609045081
$ws.Range('D28').Value = 1.068997363417893
$ws.Range('E28').Value = 0.382621474409821
$ws.Range('F28').Value = 3
$ws.Range('G28').Value = 'adaptive\_filter: wRLS, error\_metric: CPPM, fuzzy\_operator: prod, num\_generations: 5, num\_parents\_mating: 5, parallel\_processing: 10, rules: 3, sol\_per\_pop: 5 \\'

# Row 29
$ws.Range('C29').Value = 0.2200974688435089
$ws.Range('D29').Value = 0.9783457952016263
$ws.Range('E29').Value = 0.3679545340620224
$ws.Range('G29').Value = 'combination: median, n\_estimators: 50 \\'

# Row 30
$ws.Range('C30').Value = 0.2078967309800335
$ws.Range('D30').Value = 0.9241128199211365
$ws.Range('E30').Value = 0.3579224214918503

# Row 31
$ws.Range('C31').Value = 0.2101438001793833
$ws.Range('D31').Value = 0.934101170601689
$ws.Range('E31').Value = 0.3617992240005157
